$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

$xml = $xml.Replace('<Relationship Id="rId9" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/ckreibich/scholar.py" TargetMode="External"/>', '<Relationship Id="@@rId12@@" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/ckreibich/scholar.py" TargetMode="External"/>') # rels-elem-rId9
$xml = $xml.Replace('<Relationship Id="rId10" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/gimoya/theBioBucket-Archives/blob/master/R/Functions/GScholarScraper_3.1.R" TargetMode="External"/>', '<Relationship Id="@@rId13@@" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/gimoya/theBioBucket-Archives/blob/master/R/Functions/GScholarScraper_3.1.R" TargetMode="External"/>') # rels-elem-rId10
$xml = $xml.Replace('<Relationship Id="rId11" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.nihlibrary.nih.gov/services/systematic-reviews/resources" TargetMode="External"/>', '<Relationship Id="@@rId14@@" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.nihlibrary.nih.gov/services/systematic-reviews/resources" TargetMode="External"/>') # rels-elem-rId11
$xml = $xml.Replace('<Relationship Id="rId12" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.nihlibrary.nih.gov/services/systematic-reviews/systematic-review-standards-organizations" TargetMode="External"/>', '<Relationship Id="@@rId15@@" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://www.nihlibrary.nih.gov/services/systematic-reviews/systematic-review-standards-organizations" TargetMode="External"/>') # rels-elem-rId12
$xml = $xml.Replace('<Relationship Id="rId13" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/fontTable" Target="fontTable.xml"/>', '<Relationship Id="@@rId16@@" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/fontTable" Target="fontTable.xml"/>') # rels-elem-rId13
$xml = $xml.Replace('<Relationship Id="rId14" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/theme" Target="theme/theme1.xml"/>', '<Relationship Id="@@rId17@@" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/theme" Target="theme/theme1.xml"/>') # rels-elem-rId14
$xml = $xml.Replace('<w:hyperlink r:id="rId9" w:history="1">', '<w:hyperlink r:id="@@rId12@@" w:history="1">') # hyperlink-ref-rId9
$xml = $xml.Replace('<w:hyperlink r:id="rId10" w:history="1">', '<w:hyperlink r:id="@@rId13@@" w:history="1">') # hyperlink-ref-rId10
$xml = $xml.Replace('<w:hyperlink r:id="rId11" w:history="1">', '<w:hyperlink r:id="@@rId14@@" w:history="1">') # hyperlink-ref-rId11
$xml = $xml.Replace('<w:hyperlink r:id="rId12" w:history="1">', '<w:hyperlink r:id="@@rId15@@" w:history="1">') # hyperlink-ref-rId12
$xml = $xml.Replace('@@rId12@@', 'rId12') # resolve-rId12
$xml = $xml.Replace('@@rId13@@', 'rId13') # resolve-rId13
$xml = $xml.Replace('@@rId14@@', 'rId14') # resolve-rId14
$xml = $xml.Replace('@@rId15@@', 'rId15') # resolve-rId15
$xml = $xml.Replace('@@rId16@@', 'rId16') # resolve-rId16
$xml = $xml.Replace('@@rId17@@', 'rId17') # resolve-rId17
$xml = $xml.Replace('</Relationships>', '<Relationship Id="rId9" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/ropensci/RSelenium" TargetMode="External"/><Relationship Id="rId10" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://stackoverflow.com/questions/27754051/how-to-retrieve-informations-about-journals-from-isi-web-of-knowledge" TargetMode="External"/><Relationship Id="rId11" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/ropensci/webservices" TargetMode="External"/></Relationships>') # add-new-rels
$xml = $xml.Replace('<w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3570" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="3BB3A045"', '<w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3570" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="3BB3A045"') # remove-old-gobak
$xml = $xml.Replace('<w:tr w:rsidR="00177E20" w14:paraId="4C33C8BE" w14:textId="77777777" w:rsidTr="00AA48FA"><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:trHeight w:val="576"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="2610" w:type="dxa"/><w:tcBorders><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="4121033D" w14:textId="77777777" w:rsidR="00177E20" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:r><w:t>Web of Science</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4950" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="2A860AD0" w14:textId="77777777" w:rsidR="00177E20" w:rsidRPr="000260BE" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r w:rsidRPr="00177E20"><w:t>https://github.com/kousu/isi</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3570" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="73598442" w14:textId="77777777" w:rsidR="00177E20" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr></w:p></w:tc></w:tr>', '<w:tr w:rsidR="00177E20" w14:paraId="4C33C8BE" w14:textId="77777777" w:rsidTr="00AA48FA"><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:trHeight w:val="576"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="2610" w:type="dxa"/><w:tcBorders><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="4121033D" w14:textId="77777777" w:rsidR="00177E20" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:r><w:t>Web of Science</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4950" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="2A860AD0" w14:textId="77777777" w:rsidR="00177E20" w:rsidRPr="000260BE" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:hyperlink r:id="rId9" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://github.com/ropensci/RSelenium</w:t></w:r></w:hyperlink></w:p><w:p w14:paraId="2A860AD1" w14:textId="77777777" w:rsidR="00177E20" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:hyperlink r:id="rId10" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://stackoverflow.com/questions/27754051/how-to-retrieve-informations-about-journals-from-isi-web-of-knowledge</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p w14:paraId="2A860AD2" w14:textId="77777777" w:rsidR="00177E20" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:hyperlink r:id="rId11" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://github.com/ropensci/webservices</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3570" w:type="dxa"/><w:tcBorders><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p w14:paraId="73598442" w14:textId="77777777" w:rsidR="00177E20" w:rsidRDefault="00177E20" w:rsidP="00A61AAB"><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>ResearchID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:b/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>V-7730-2018</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:tc></w:tr>') # web-of-science-row
$xml = $xml.Replace('<w:tblGrid><w:gridCol w:w="2486"/><w:gridCol w:w="5310"/><w:gridCol w:w="3334"/></w:tblGrid>', '<w:tblGrid><w:gridCol w:w="2482"/><w:gridCol w:w="5310"/><w:gridCol w:w="3338"/></w:tblGrid>') # tblgrid-widths

$d.Content.InsertXML($xml)
Write-Output "applied edits"
